$wb = $excel.ActiveWorkbook

# --- Sheet scaffolding: rename/reorder to the 4-sheet layout ----------
# NOTE: worksheet handles returned by Worksheets.Item(...) can become stale
# positional references once Worksheets.Add() shifts sheets around, so every
# handle below is (re-)fetched by name right before it is used.
$wsBattingOld = $wb.Worksheets.Item("ODI Batting")

# New sheet inserted in front of 'ODI Batting' -> becomes sheet 1
$wsPlayerInfo = $wb.Worksheets.Add($wsBattingOld)
$wsPlayerInfo.Name = "Player Info"

# New sheet appended after the last sheet -> becomes the last sheet
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add($null, $wsLast)
$wsExtra.Name = "ODI Batting Extra"

# Re-fetch the pre-existing sheets fresh, by name, now that the sheet
# collection has settled into its final shape/order.
$wsPlayerInfo = $wb.Worksheets.Item("Player Info")
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsExtra = $wb.Worksheets.Item("ODI Batting Extra")

# --- Player Info (new sheet) -----------------------------------------------------
$wsPlayerInfo.Range("A1").NumberFormat = "@"
$wsPlayerInfo.Range("A1").Value = "ID"
$wsPlayerInfo.Range("B1").NumberFormat = "@"
$wsPlayerInfo.Range("B1").Value = "NAME"
$wsPlayerInfo.Range("C1").NumberFormat = "@"
$wsPlayerInfo.Range("C1").Value = "BATTING_HAND"
$wsPlayerInfo.Range("D1").NumberFormat = "@"
$wsPlayerInfo.Range("D1").Value = "BOWL_STYLE"
$wsPlayerInfo.Range("A2").NumberFormat = "@"
$wsPlayerInfo.Range("A2").Value = "5935"
$wsPlayerInfo.Range("B2").NumberFormat = "@"
$wsPlayerInfo.Range("B2").Value = "Akeal Jerome Hosein"
$wsPlayerInfo.Range("C2").NumberFormat = "@"
$wsPlayerInfo.Range("C2").Value = "Left Handed"
$wsPlayerInfo.Range("D2").NumberFormat = "@"
$wsPlayerInfo.Range("D2").Value = "Left Arm Orthodox"

# --- ODI Batting: cell updates only ----------------------------------
$wsBatting.Range("D1").NumberFormat = "@"
$wsBatting.Range("D1").Value = "MATCH_CODE"
$wsBatting.Range("D2").NumberFormat = "@"
$wsBatting.Range("D2").Value = "4443"
$wsBatting.Range("D3").NumberFormat = "@"
$wsBatting.Range("D3").Value = "4445"
$wsBatting.Range("D4").NumberFormat = "@"
$wsBatting.Range("D4").Value = "4447"
$wsBatting.Range("D5").NumberFormat = "@"
$wsBatting.Range("D5").Value = "4449"
$wsBatting.Range("D6").NumberFormat = "@"
$wsBatting.Range("D6").Value = "4450"
$wsBatting.Range("D7").NumberFormat = "@"
$wsBatting.Range("D7").Value = "4451"
$wsBatting.Range("D8").NumberFormat = "@"
$wsBatting.Range("D8").Value = "4483"
$wsBatting.Range("D9").NumberFormat = "@"
$wsBatting.Range("D9").Value = "4484"
$wsBatting.Range("D10").NumberFormat = "@"
$wsBatting.Range("D10").Value = "4486"
$wsBatting.Range("D11").NumberFormat = "@"
$wsBatting.Range("D11").Value = "4519"
$wsBatting.Range("D12").NumberFormat = "@"
$wsBatting.Range("D12").Value = "4520"
$wsBatting.Range("D13").NumberFormat = "@"
$wsBatting.Range("D13").Value = "4522"
$wsBatting.Range("D14").NumberFormat = "@"
$wsBatting.Range("D14").Value = "4533"
$wsBatting.Range("D15").NumberFormat = "@"
$wsBatting.Range("D15").Value = "4535"
$wsBatting.Range("D16").NumberFormat = "@"
$wsBatting.Range("D16").Value = "4577"
$wsBatting.Range("D17").NumberFormat = "@"
$wsBatting.Range("D17").Value = "4580"
$wsBatting.Range("D18").NumberFormat = "@"
$wsBatting.Range("D18").Value = "4583"
$wsBatting.Range("D19").NumberFormat = "@"
$wsBatting.Range("D19").Value = "4586"
$wsBatting.Range("D20").NumberFormat = "@"
$wsBatting.Range("D20").Value = "4590"
$wsBatting.Range("D21").NumberFormat = "@"
$wsBatting.Range("D21").Value = "4592"
$wsBatting.Range("D22").NumberFormat = "@"
$wsBatting.Range("D22").Value = "4606"
$wsBatting.Range("D23").NumberFormat = "@"
$wsBatting.Range("D23").Value = "4611"
$wsBatting.Range("D24").NumberFormat = "@"
$wsBatting.Range("D24").Value = "4616"
$wsBatting.Range("D25").NumberFormat = "@"
$wsBatting.Range("D25").Value = "4621"
$wsBatting.Range("D26").NumberFormat = "@"
$wsBatting.Range("D26").Value = "4623"
$wsBatting.Range("D27").NumberFormat = "@"
$wsBatting.Range("D27").Value = "4624"
$wsBatting.Range("D28").NumberFormat = "@"
$wsBatting.Range("D28").Value = "4636"
$wsBatting.Range("D29").NumberFormat = "@"
$wsBatting.Range("D29").Value = "4639"
$wsBatting.Range("D30").NumberFormat = "@"
$wsBatting.Range("D30").Value = "4642"
$wsBatting.Range("D31").NumberFormat = "@"
$wsBatting.Range("D31").Value = "4727"
$wsBatting.Range("D32").NumberFormat = "@"
$wsBatting.Range("D32").Value = "4731"

# --- ODI Bowling: cell updates only ----------------------------------
$wsBowling.Range("B1").NumberFormat = "@"
$wsBowling.Range("B1").Value = "MATCH_CODE"
$wsBowling.Range("B2").NumberFormat = "@"
$wsBowling.Range("B2").Value = "4443"
$wsBowling.Range("B3").NumberFormat = "@"
$wsBowling.Range("B3").Value = "4445"
$wsBowling.Range("B4").NumberFormat = "@"
$wsBowling.Range("B4").Value = "4447"
$wsBowling.Range("B5").NumberFormat = "@"
$wsBowling.Range("B5").Value = "4449"
$wsBowling.Range("B6").NumberFormat = "@"
$wsBowling.Range("B6").Value = "4450"
$wsBowling.Range("B7").NumberFormat = "@"
$wsBowling.Range("B7").Value = "4451"
$wsBowling.Range("B8").NumberFormat = "@"
$wsBowling.Range("B8").Value = "4483"
$wsBowling.Range("B9").NumberFormat = "@"
$wsBowling.Range("B9").Value = "4484"
$wsBowling.Range("B10").NumberFormat = "@"
$wsBowling.Range("B10").Value = "4486"
$wsBowling.Range("B11").NumberFormat = "@"
$wsBowling.Range("B11").Value = "4519"
$wsBowling.Range("B12").NumberFormat = "@"
$wsBowling.Range("B12").Value = "4520"
$wsBowling.Range("B13").NumberFormat = "@"
$wsBowling.Range("B13").Value = "4522"
$wsBowling.Range("B14").NumberFormat = "@"
$wsBowling.Range("B14").Value = "4533"
$wsBowling.Range("B15").NumberFormat = "@"
$wsBowling.Range("B15").Value = "4535"
$wsBowling.Range("B16").NumberFormat = "@"
$wsBowling.Range("B16").Value = "4577"
$wsBowling.Range("B17").NumberFormat = "@"
$wsBowling.Range("B17").Value = "4580"
$wsBowling.Range("B18").NumberFormat = "@"
$wsBowling.Range("B18").Value = "4583"
$wsBowling.Range("B19").NumberFormat = "@"
$wsBowling.Range("B19").Value = "4586"
$wsBowling.Range("B20").NumberFormat = "@"
$wsBowling.Range("B20").Value = "4590"
$wsBowling.Range("B21").NumberFormat = "@"
$wsBowling.Range("B21").Value = "4592"
$wsBowling.Range("B22").NumberFormat = "@"
$wsBowling.Range("B22").Value = "4606"
$wsBowling.Range("B23").NumberFormat = "@"
$wsBowling.Range("B23").Value = "4611"
$wsBowling.Range("B24").NumberFormat = "@"
$wsBowling.Range("B24").Value = "4616"
$wsBowling.Range("B25").NumberFormat = "@"
$wsBowling.Range("B25").Value = "4621"
$wsBowling.Range("B26").NumberFormat = "@"
$wsBowling.Range("B26").Value = "4623"
$wsBowling.Range("B27").NumberFormat = "@"
$wsBowling.Range("B27").Value = "4624"
$wsBowling.Range("B28").NumberFormat = "@"
$wsBowling.Range("B28").Value = "4636"
$wsBowling.Range("B29").NumberFormat = "@"
$wsBowling.Range("B29").Value = "4639"
$wsBowling.Range("B30").NumberFormat = "@"
$wsBowling.Range("B30").Value = "4642"
$wsBowling.Range("B31").NumberFormat = "@"
$wsBowling.Range("B31").Value = "4727"
$wsBowling.Range("B32").NumberFormat = "@"
$wsBowling.Range("B32").Value = "4731"

# --- ODI Batting Extra (new sheet) -----------------------------------------------------
$wsExtra.Range("A1").NumberFormat = "@"
$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").NumberFormat = "@"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").NumberFormat = "@"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").NumberFormat = "@"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").NumberFormat = "@"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").NumberFormat = "@"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"
$wsExtra.Range("A2").NumberFormat = "@"
$wsExtra.Range("A2").Value = "4522"
$wsExtra.Range("B2").Value = ""
$wsExtra.Range("C2").Value = ""
$wsExtra.Range("D2").Value = ""
$wsExtra.Range("E2").Value = ""
$wsExtra.Range("F2").NumberFormat = "@"
$wsExtra.Range("F2").Value = "NO"
$wsExtra.Range("A3").NumberFormat = "@"
$wsExtra.Range("A3").Value = "4533"
$wsExtra.Range("B3").Value = ""
$wsExtra.Range("C3").Value = ""
$wsExtra.Range("D3").Value = ""
$wsExtra.Range("E3").Value = ""
$wsExtra.Range("F3").NumberFormat = "@"
$wsExtra.Range("F3").Value = "NO"
$wsExtra.Range("A4").NumberFormat = "@"
$wsExtra.Range("A4").Value = "4535"
$wsExtra.Range("B4").Value = 7
$wsExtra.Range("C4").NumberFormat = "@"
$wsExtra.Range("C4").Value = "3"
$wsExtra.Range("D4").NumberFormat = "@"
$wsExtra.Range("D4").Value = "0"
$wsExtra.Range("E4").NumberFormat = "@"
$wsExtra.Range("E4").Value = "17.62%"
$wsExtra.Range("F4").NumberFormat = "@"
$wsExtra.Range("F4").Value = "NO"
$wsExtra.Range("A5").NumberFormat = "@"
$wsExtra.Range("A5").Value = "4577"
$wsExtra.Range("B5").Value = 6
$wsExtra.Range("C5").NumberFormat = "@"
$wsExtra.Range("C5").Value = "0"
$wsExtra.Range("D5").NumberFormat = "@"
$wsExtra.Range("D5").Value = "0"
$wsExtra.Range("E5").NumberFormat = "@"
$wsExtra.Range("E5").Value = "0.97%"
$wsExtra.Range("F5").NumberFormat = "@"
$wsExtra.Range("F5").Value = "NO"
$wsExtra.Range("A6").NumberFormat = "@"
$wsExtra.Range("A6").Value = "4580"
$wsExtra.Range("B6").Value = ""
$wsExtra.Range("C6").Value = ""
$wsExtra.Range("D6").Value = ""
$wsExtra.Range("E6").Value = ""
$wsExtra.Range("F6").NumberFormat = "@"
$wsExtra.Range("F6").Value = "NO"
$wsExtra.Range("A7").NumberFormat = "@"
$wsExtra.Range("A7").Value = "4583"
$wsExtra.Range("B7").Value = ""
$wsExtra.Range("C7").Value = ""
$wsExtra.Range("D7").Value = ""
$wsExtra.Range("E7").Value = ""
$wsExtra.Range("F7").NumberFormat = "@"
$wsExtra.Range("F7").Value = "NO"
$wsExtra.Range("A8").NumberFormat = "@"
$wsExtra.Range("A8").Value = "4586"
$wsExtra.Range("B8").Value = 7
$wsExtra.Range("C8").NumberFormat = "@"
$wsExtra.Range("C8").Value = "2"
$wsExtra.Range("D8").NumberFormat = "@"
$wsExtra.Range("D8").Value = "6"
$wsExtra.Range("E8").NumberFormat = "@"
$wsExtra.Range("E8").Value = "27.78%"
$wsExtra.Range("F8").NumberFormat = "@"
$wsExtra.Range("F8").Value = "NO"
$wsExtra.Range("A9").NumberFormat = "@"
$wsExtra.Range("A9").Value = "4590"
$wsExtra.Range("B9").Value = ""
$wsExtra.Range("C9").Value = ""
$wsExtra.Range("D9").Value = ""
$wsExtra.Range("E9").Value = ""
$wsExtra.Range("F9").NumberFormat = "@"
$wsExtra.Range("F9").Value = "NO"
$wsExtra.Range("A10").NumberFormat = "@"
$wsExtra.Range("A10").Value = "4592"
$wsExtra.Range("B10").Value = 7
$wsExtra.Range("C10").NumberFormat = "@"
$wsExtra.Range("C10").Value = "0"
$wsExtra.Range("D10").NumberFormat = "@"
$wsExtra.Range("D10").Value = "0"
$wsExtra.Range("E10").NumberFormat = "@"
$wsExtra.Range("E10").Value = "1.85%"
$wsExtra.Range("F10").NumberFormat = "@"
$wsExtra.Range("F10").Value = "NO"
$wsExtra.Range("A11").NumberFormat = "@"
$wsExtra.Range("A11").Value = "4606"
$wsExtra.Range("B11").Value = 8
$wsExtra.Range("C11").NumberFormat = "@"
$wsExtra.Range("C11").Value = "0"
$wsExtra.Range("D11").NumberFormat = "@"
$wsExtra.Range("D11").Value = "0"
$wsExtra.Range("E11").NumberFormat = "@"
$wsExtra.Range("E11").Value = "0.56%"
$wsExtra.Range("F11").NumberFormat = "@"
$wsExtra.Range("F11").Value = "NO"
$wsExtra.Range("A12").NumberFormat = "@"
$wsExtra.Range("A12").Value = "4611"
$wsExtra.Range("B12").Value = 7
$wsExtra.Range("C12").NumberFormat = "@"
$wsExtra.Range("C12").Value = "2"
$wsExtra.Range("D12").NumberFormat = "@"
$wsExtra.Range("D12").Value = "0"
$wsExtra.Range("E12").NumberFormat = "@"
$wsExtra.Range("E12").Value = "10.49%"
$wsExtra.Range("F12").NumberFormat = "@"
$wsExtra.Range("F12").Value = "NO"
$wsExtra.Range("A13").NumberFormat = "@"
$wsExtra.Range("A13").Value = "4616"
$wsExtra.Range("B13").Value = ""
$wsExtra.Range("C13").Value = ""
$wsExtra.Range("D13").Value = ""
$wsExtra.Range("E13").Value = ""
$wsExtra.Range("F13").NumberFormat = "@"
$wsExtra.Range("F13").Value = "NO"
$wsExtra.Range("A14").NumberFormat = "@"
$wsExtra.Range("A14").Value = "4621"
$wsExtra.Range("B14").Value = 8
$wsExtra.Range("C14").NumberFormat = "@"
$wsExtra.Range("C14").Value = "0"
$wsExtra.Range("D14").NumberFormat = "@"
$wsExtra.Range("D14").Value = "0"
$wsExtra.Range("E14").NumberFormat = "@"
$wsExtra.Range("E14").Value = "0.73%"
$wsExtra.Range("F14").NumberFormat = "@"
$wsExtra.Range("F14").Value = "NO"
$wsExtra.Range("A15").NumberFormat = "@"
$wsExtra.Range("A15").Value = "4623"
$wsExtra.Range("B15").Value = 9
$wsExtra.Range("C15").Value = ""
$wsExtra.Range("D15").Value = ""
$wsExtra.Range("E15").Value = ""
$wsExtra.Range("F15").NumberFormat = "@"
$wsExtra.Range("F15").Value = "NO"
$wsExtra.Range("A16").NumberFormat = "@"
$wsExtra.Range("A16").Value = "4624"
$wsExtra.Range("B16").Value = 9
$wsExtra.Range("C16").NumberFormat = "@"
$wsExtra.Range("C16").Value = "0"
$wsExtra.Range("D16").NumberFormat = "@"
$wsExtra.Range("D16").Value = "1"
$wsExtra.Range("E16").NumberFormat = "@"
$wsExtra.Range("E16").Value = "6.83%"
$wsExtra.Range("F16").NumberFormat = "@"
$wsExtra.Range("F16").Value = "NO"
$wsExtra.Range("A17").NumberFormat = "@"
$wsExtra.Range("A17").Value = "4636"
$wsExtra.Range("B17").Value = ""
$wsExtra.Range("C17").Value = ""
$wsExtra.Range("D17").Value = ""
$wsExtra.Range("E17").Value = ""
$wsExtra.Range("F17").NumberFormat = "@"
$wsExtra.Range("F17").Value = "NO"
$wsExtra.Range("A18").NumberFormat = "@"
$wsExtra.Range("A18").Value = "4639"
$wsExtra.Range("B18").Value = ""
$wsExtra.Range("C18").Value = ""
$wsExtra.Range("D18").Value = ""
$wsExtra.Range("E18").Value = ""
$wsExtra.Range("F18").NumberFormat = "@"
$wsExtra.Range("F18").Value = "NO"
$wsExtra.Range("A19").NumberFormat = "@"
$wsExtra.Range("A19").Value = "4642"
$wsExtra.Range("B19").Value = 9
$wsExtra.Range("C19").NumberFormat = "@"
$wsExtra.Range("C19").Value = "1"
$wsExtra.Range("D19").NumberFormat = "@"
$wsExtra.Range("D19").Value = "0"
$wsExtra.Range("E19").NumberFormat = "@"
$wsExtra.Range("E19").Value = "5.38%"
$wsExtra.Range("F19").NumberFormat = "@"
$wsExtra.Range("F19").Value = "NO"
$wsExtra.Range("A20").NumberFormat = "@"
$wsExtra.Range("A20").Value = "4727"
$wsExtra.Range("B20").Value = ""
$wsExtra.Range("C20").Value = ""
$wsExtra.Range("D20").Value = ""
$wsExtra.Range("E20").Value = ""
$wsExtra.Range("F20").Value = ""
$wsExtra.Range("A21").NumberFormat = "@"
$wsExtra.Range("A21").Value = "4731"
$wsExtra.Range("B21").Value = ""
$wsExtra.Range("C21").Value = ""
$wsExtra.Range("D21").Value = ""
$wsExtra.Range("E21").Value = ""
$wsExtra.Range("F21").Value = ""

# --- Header styling: bold + centered, matching the workbook's house style
$wsPlayerInfo.Range("A1:D1").Font.Bold = $true
$wsPlayerInfo.Range("A1:D1").HorizontalAlignment = -4108
$wsPlayerInfo.Range("A1:D1").VerticalAlignment = -4160
$wsPlayerInfo.Range("A1:D1").Borders.LineStyle = 1
$wsExtra.Range("A1:F1").Font.Bold = $true
$wsExtra.Range("A1:F1").HorizontalAlignment = -4108
$wsExtra.Range("A1:F1").VerticalAlignment = -4160
$wsExtra.Range("A1:F1").Borders.LineStyle = 1

$wsPlayerInfo.Range("A1").Select()
